# "V. 71 'Mi año en Oxford'"
#
# Adds a new rated movie ("Mi año en Oxford") to the "Películas" sheet's
# table (Tabla24, sorted descending by the average-score column C).
#
# The previously most-recently-added entry, "Chappie" (row 67), carried a
# slightly different cell style on its title cell (a leftover artifact from
# when *it* was the newest addition). That marker style normally moves on to
# whichever title cell is typed next, so we replicate that: the new row's
# title cell inherits the special style and "Chappie" gets normalized back
# to the regular title style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")
$ws.Activate()

# Insert a fresh row right after "Chappie" (row 67) so the new row copies
# that row's formatting (including the special title-cell style) before we
# let the table's sort relocate it to its proper ranked position.
$ws.Rows("68:68").Insert(-4121, 0)  # xlShiftDown, xlFormatFromLeftOrAbove

# Fill in the new movie's data.
$ws.Range("B68").Value = "Mi año en Oxford"
$ws.Range("C68").Formula = "=AVERAGE(D68,E68,E68,F68,G68,H68,H68,I68)"
$ws.Range("D68").Value = 6
$ws.Range("E68").Value = 3
$ws.Range("F68").Value = 5
$ws.Range("G68").Value = 6
$ws.Range("H68").Value = 6
$ws.Range("I68").Value = 4.9

# Normalize "Chappie"'s title-cell style back to the regular one now that it
# is no longer the newest entry (re-enter its value to drop the old style).
$chappie = $ws.Range("B67")
$chappieTitle = $chappie.Value()
$chappie.Cut($chappie)
$chappie.Value = $chappieTitle

# Grow the table to cover the newly inserted row.
$lo = $ws.ListObjects.Item("Tabla24")
$lo.Resize($ws.Range("B2:I95"))

# Re-apply the table's existing sort (by total score, descending) so the new
# row slots into its correct ranked position and everything below shifts
# down.
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("C2:C95"), 0, 2) | Out-Null
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Match the author's final selection/view state.
$ws.Range("B94").Select()
